$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking strings
# (e.g. "313.12", "1.00") are stored as text, matching the source
# workbook where these are inline strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '42.608.95'
$ws.Cells.Item(2, 5).Value = '  -0.41%  '
$ws.Cells.Item(3, 4).Value = '2.544.83'
$ws.Cells.Item(3, 5).Value = '  +0.71%  '
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
$ws.Cells.Item(5, 4).Value = '313.12'
$ws.Cells.Item(5, 5).Value = '  -1.01%  '
$ws.Cells.Item(6, 4).Value = '101.17'
$ws.Cells.Item(6, 5).Value = '  +5.66%  '
$ws.Cells.Item(7, 5).Value = '  -0.82%  '
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).Value = '0.529'
$ws.Cells.Item(9, 5).Value = '  -1.62%  '
$ws.Cells.Item(10, 4).Value = '36.04'
$ws.Cells.Item(10, 5).Value = '  +1.80%  '
$ws.Cells.Item(11, 4).Value = '0.0804'
$ws.Cells.Item(11, 5).Value = '  -0.64%  '
$ws.Cells.Item(12, 5).Value = '  -0.89%  '
$ws.Cells.Item(13, 5).Value = '  -0.32%  '
$ws.Cells.Item(14, 4).Value = '2.935.04'
$ws.Cells.Item(14, 5).Value = '  +0.38%  '
$ws.Cells.Item(15, 4).Value = '15.92'
$ws.Cells.Item(15, 5).Value = '  +6.50%  '
$ws.Cells.Item(16, 4).Value = '2.570.64'
$ws.Cells.Item(16, 5).Value = '  +3.45%  '
$ws.Cells.Item(17, 4).Value = '0.833'
$ws.Cells.Item(17, 5).Value = '  -1.49%  '
$ws.Cells.Item(18, 4).Value = '42.624.69'
$ws.Cells.Item(18, 5).Value = '  -0.64%  '
$ws.Cells.Item(19, 5).Value = '  +0.42%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0954'
$ws.Cells.Item(20, 5).Value = '  -0.80%  '
$ws.Cells.Item(21, 5).Value = '  -1.43%  '
$ws.Cells.Item(22, 4).Value = '69.09'
$ws.Cells.Item(22, 5).Value = '  -0.53%  '
$ws.Cells.Item(23, 4).Value = '244.61'
$ws.Cells.Item(23, 5).Value = '  -3.21%  '
$ws.Cells.Item(24, 4).Value = '2.93'
$ws.Cells.Item(24, 5).Value = '  -0.40%  '
$ws.Cells.Item(25, 5).Value = '  +0.88%  '
$ws.Cells.Item(26, 4).Value = '26.40'
$ws.Cells.Item(26, 5).Value = '  -2.11%  '
$ws.Cells.Item(27, 5).Value = '  +0.02%  '
$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).Value = '2.35'
$ws.Cells.Item(28, 5).Value = '  -1.63%  '
$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).Value = '40.43'
$ws.Cells.Item(29, 5).Value = '  +0.05%  '
$ws.Cells.Item(30, 4).Value = '10.10'
$ws.Cells.Item(30, 5).Value = '  -1.72%  '
$ws.Cells.Item(31, 4).Value = '158.52'
$ws.Cells.Item(31, 5).Value = '  +1.59%  '
$ws.Cells.Item(32, 5).Value = '  -2.54%  '
$ws.Cells.Item(33, 4).Value = '2.79'
$ws.Cells.Item(33, 5).Value = '  +13.94%  '
$ws.Cells.Item(34, 5).Value = '  +1.27%  '
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).Value = '2.06'
$ws.Cells.Item(35, 5).Value = '  -1.27%  '
$ws.Cells.Item(36, 2).Value = 'WEMIXToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(36, 4).Value = '2.62'
$ws.Cells.Item(36, 5).Value = '  -3.02%  '
$ws.Cells.Item(37, 4).Value = '3.22'
$ws.Cells.Item(37, 5).Value = '  -3.48%  '
$ws.Cells.Item(38, 4).Value = '18.14'
$ws.Cells.Item(38, 5).Value = '  -5.14%  '
$ws.Cells.Item(39, 5).Value = '  -1.53%  '
$ws.Cells.Item(40, 5).Value = '  -0.32%  '
$ws.Cells.Item(41, 4).Value = '22.26'
$ws.Cells.Item(41, 5).Value = '  +1.95%  '
$ws.Cells.Item(42, 4).Value = '4.14'
$ws.Cells.Item(42, 5).Value = '  +9.01%  '
$ws.Cells.Item(43, 5).Value = '  +0.05%  '
$ws.Cells.Item(44, 4).Value = '3.32'
$ws.Cells.Item(44, 5).Value = '  +2.03%  '
$ws.Cells.Item(45, 4).Value = '0.0298'
$ws.Cells.Item(45, 5).Value = '  -1.32%  '
$ws.Cells.Item(46, 4).Value = '1.970.94'
$ws.Cells.Item(46, 5).Value = '  -1.22%  '
$ws.Cells.Item(47, 4).Value = '8.88'
$ws.Cells.Item(47, 5).Value = '  -1.71%  '
$ws.Cells.Item(48, 4).Value = '2.793.59'
$ws.Cells.Item(48, 5).Value = '  +0.56%  '
$ws.Cells.Item(49, 4).Value = '81.08'
$ws.Cells.Item(49, 5).Value = '  -3.88%  '
$ws.Cells.Item(50, 5).Value = '  +0.94%  '
$ws.Cells.Item(51, 4).Value = '73.30'
$ws.Cells.Item(51, 5).Value = '  -0.98%  '

# Restore the original (default/"Normal") style on column D so no
# stray style index is left on these cells.
$ws.Range("D2:D51").Style = "Normal"
